$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27, column A: was stored as text "0.01" -> becomes a genuine number 0.01
$ws.Range("A27").Value = 0.01

# New row 28 (all plain numbers)
$ws.Range("A28").Value = 0.01
$ws.Range("B28").Value = 0.08400000000000001
$ws.Range("C28").Value = 18431
$ws.Range("D28").Value = -0.86

# New row 29 - column A must stay textual ("0.06"), the rest are numbers
$a29 = $ws.Range("A29")
$a29.NumberFormat = "@"
$a29.Value = "0.06"
$a29.NumberFormat = "General"
$a29.Style = "Normal"

$ws.Range("B29").Value = 0.031
$ws.Range("C29").Value = 17528
$ws.Range("D29").Value = -4.9
